$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 8.280371333333333
$ws.Range("H2").Value = 24.841114
$ws.Range("I2").Value = 0.2946400644635011
$ws.Range("J2").Value = 0.3116548779253407
$ws.Range("M2").Value = 44.04223000000001
$ws.Range("N2").Value = 132.12669
$ws.Range("O2").Value = 0.1792438957545786
$ws.Range("P2").Value = 0.1868246872369915
$ws.Range("Q2").Value = 364.6860187480734
$ws.Range("R2").Value = 3282.174168732661
$ws.Range("S2").Value = 0.05281243299981812
$ws.Range("T2").Value = 0.05822482509428456
$ws.Range("G3").Value = 8.280371333333333
$ws.Range("H3").Value = 24.841114
$ws.Range("I3").Value = 0.2946400644635011
$ws.Range("J3").Value = 0.3116548779253407
$ws.Range("O3").Value = 0.2018201397722426
$ws.Range("P3").Value = 0.2103557520458098
$ws.Range("Q3").Value = 410.6191899415842
$ws.Range("R3").Value = 3695.572709474258
$ws.Range("S3").Value = 0.05946429899252635
$ws.Range("T3").Value = 0.0655583962247301
$ws.Range("G4").Value = 8.280371333333333
$ws.Range("H4").Value = 24.841114
$ws.Range("I4").Value = 0.2946400644635011
$ws.Range("J4").Value = 0.3116548779253407
$ws.Range("M4").Value = 51.56497066666667
$ws.Range("N4").Value = 154.694912
$ws.Range("O4").Value = 0.2098600871655206
$ws.Range("P4").Value = 0.2187357342528896
$ws.Range("Q4").Value = 426.977104912441
$ws.Range("R4").Value = 3842.793944211969
$ws.Range("S4").Value = 0.06183318961076494
$ws.Range("T4").Value = 0.06817005855649408
$ws.Range("G5").Value = 8.280371333333333
$ws.Range("H5").Value = 24.841114
$ws.Range("I5").Value = 0.2946400644635011
$ws.Range("J5").Value = 0.3116548779253407
$ws.Range("M5").Value = 29.9106925
$ws.Range("N5").Value = 59.821385
$ws.Range("O5").Value = 0.1217310987299521
$ws.Range("P5").Value = 0.08458632803643724
$ws.Range("Q5").Value = 247.6716407371483
$ws.Range("R5").Value = 1486.02984442289
$ws.Range("S5").Value = 0.03586685877700589
$ws.Range("T5").Value = 0.02636174173834868
$ws.Range("G6").Value = 8.280371333333333
$ws.Range("H6").Value = 24.841114
$ws.Range("I6").Value = 0.2946400644635011
$ws.Range("J6").Value = 0.3116548779253407
$ws.Range("M6").Value = 70.603826
$ws.Range("N6").Value = 211.811478
$ws.Range("O6").Value = 0.2873447785777061
$ws.Range("P6").Value = 0.2994974984278718
$ws.Range("Q6").Value = 584.6258968340546
$ws.Range("R6").Value = 5261.633071506492
$ws.Range("S6").Value = 0.08466328408338579
$ws.Range("T6").Value = 0.09333985631148331
$ws.Range("I7").Value = 0.405746032520008
$ws.Range("J7").Value = 0.4291769704298953
$ws.Range("M7").Value = 44.04223000000001
$ws.Range("N7").Value = 132.12669
$ws.Range("O7").Value = 0.1792438957545786
$ws.Range("P7").Value = 0.1868246872369915
$ws.Range("Q7").Value = 502.2056504500867
$ws.Range("R7").Value = 4519.85085405078
$ws.Range("S7").Value = 0.07272749955585019
$ws.Range("T7").Value = 0.08018085326988475
$ws.Range("I8").Value = 0.405746032520008
$ws.Range("J8").Value = 0.4291769704298953
$ws.Range("O8").Value = 0.2018201397722426
$ws.Range("P8").Value = 0.2103557520458098
$ws.Range("S8").Value = 0.08188772099522089
$ws.Range("T8").Value = 0.0902798443755229
$ws.Range("I9").Value = 0.405746032520008
$ws.Range("J9").Value = 0.4291769704298953
$ws.Range("M9").Value = 51.56497066666667
$ws.Range("N9").Value = 154.694912
$ws.Range("O9").Value = 0.2098600871655206
$ws.Range("P9").Value = 0.2187357342528896
$ws.Range("Q9").Value = 587.9861131939272
$ws.Range("R9").Value = 5291.875018745344
$ws.Range("S9").Value = 0.08514989775171303
$ws.Range("T9").Value = 0.09387633975141382
$ws.Range("I10").Value = 0.405746032520008
$ws.Range("J10").Value = 0.4291769704298953
$ws.Range("M10").Value = 29.9106925
$ws.Range("N10").Value = 59.821385
$ws.Range("O10").Value = 0.1217310987299521
$ws.Range("P10").Value = 0.08458632803643724
$ws.Range("Q10").Value = 341.0662625933116
$ws.Range("R10").Value = 2046.39757555987
$ws.Range("S10").Value = 0.04939191034397943
$ws.Range("T10").Value = 0.03630250400646745
$ws.Range("I11").Value = 0.405746032520008
$ws.Range("J11").Value = 0.4291769704298953
$ws.Range("M11").Value = 70.603826
$ws.Range("N11").Value = 211.811478
$ws.Range("O11").Value = 0.2873447785777061
$ws.Range("P11").Value = 0.2994974984278718
$ws.Range("Q11").Value = 805.0827662585372
$ws.Range("R11").Value = 7245.744896326835
$ws.Range("S11").Value = 0.1165890038732445
$ws.Range("T11").Value = 0.1285374290266063
$ws.Range("G12").Value = 1.864050333333333
$ws.Range("H12").Value = 5.592150999999999
$ws.Range("I12").Value = 0.06632841551025578
$ws.Range("J12").Value = 0.07015873512134246
$ws.Range("M12").Value = 44.04223000000001
$ws.Range("N12").Value = 132.12669
$ws.Range("O12").Value = 0.1792438957545786
$ws.Range("P12").Value = 0.1868246872369915
$ws.Range("Q12").Value = 82.09693351224334
$ws.Range("R12").Value = 738.8724016101901
$ws.Range("S12").Value = 0.01188896359528666
$ws.Range("T12").Value = 0.01310738374598774
$ws.Range("G13").Value = 1.864050333333333
$ws.Range("H13").Value = 5.592150999999999
$ws.Range("I13").Value = 0.06632841551025578
$ws.Range("J13").Value = 0.07015873512134246
$ws.Range("O13").Value = 0.2018201397722426
$ws.Range("P13").Value = 0.2103557520458098
$ws.Range("Q13").Value = 92.43726000577188
$ws.Range("R13").Value = 831.9353400519469
$ws.Range("S13").Value = 0.0133864100891512
$ws.Range("T13").Value = 0.01475829348903276
$ws.Range("G14").Value = 1.864050333333333
$ws.Range("H14").Value = 5.592150999999999
$ws.Range("I14").Value = 0.06632841551025578
$ws.Range("J14").Value = 0.07015873512134246
$ws.Range("M14").Value = 51.56497066666667
$ws.Range("N14").Value = 154.694912
$ws.Range("O14").Value = 0.2098600871655206
$ws.Range("P14").Value = 0.2187357342528896
$ws.Range("Q14").Value = 96.11970075952355
$ws.Range("R14").Value = 865.077306835712
$ws.Range("S14").Value = 0.01391968706053314
$ws.Range("T14").Value = 0.01534622244102084
$ws.Range("G15").Value = 1.864050333333333
$ws.Range("H15").Value = 5.592150999999999
$ws.Range("I15").Value = 0.06632841551025578
$ws.Range("J15").Value = 0.07015873512134246
$ws.Range("M15").Value = 29.9106925
$ws.Range("N15").Value = 59.821385
$ws.Range("O15").Value = 0.1217310987299521
$ws.Range("P15").Value = 0.08458632803643724
$ws.Range("Q15").Value = 55.75503632485583
$ws.Range("R15").Value = 334.530217949135
$ws.Range("S15").Value = 0.008074230897080229
$ws.Range("T15").Value = 0.005934469783595384
$ws.Range("G16").Value = 1.864050333333333
$ws.Range("H16").Value = 5.592150999999999
$ws.Range("I16").Value = 0.06632841551025578
$ws.Range("J16").Value = 0.07015873512134246
$ws.Range("M16").Value = 70.603826
$ws.Range("N16").Value = 211.811478
$ws.Range("O16").Value = 0.2873447785777061
$ws.Range("P16").Value = 0.2994974984278718
$ws.Range("Q16").Value = 131.6090853899086
$ws.Range("R16").Value = 1184.481768509178
$ws.Range("S16").Value = 0.01905912386820454
$ws.Range("T16").Value = 0.02101236566170574
$ws.Range("G17").Value = 4.6029105
$ws.Range("H17").Value = 9.205821
$ws.Range("I17").Value = 0.1637851482553954
$ws.Range("J17").Value = 0.1154955860658076
$ws.Range("M17").Value = 44.04223000000001
$ws.Range("N17").Value = 132.12669
$ws.Range("O17").Value = 0.1792438957545786
$ws.Range("P17").Value = 0.1868246872369915
$ws.Range("Q17").Value = 202.7224429104151
$ws.Range("R17").Value = 1216.33465746249
$ws.Range("S17").Value = 0.0293574880400383
$ws.Range("T17").Value = 0.02157742674399754
$ws.Range("G18").Value = 4.6029105
$ws.Range("H18").Value = 9.205821
$ws.Range("I18").Value = 0.1637851482553954
$ws.Range("J18").Value = 0.1154955860658076
$ws.Range("O18").Value = 0.2018201397722426
$ws.Range("P18").Value = 0.2103557520458098
$ws.Range("Q18").Value = 228.2558722064895
$ws.Range("R18").Value = 1369.535233238937
$ws.Range("S18").Value = 0.03305514151352137
$ws.Range("T18").Value = 0.02429516086484451
$ws.Range("G19").Value = 4.6029105
$ws.Range("H19").Value = 9.205821
$ws.Range("I19").Value = 0.1637851482553954
$ws.Range("J19").Value = 0.1154955860658076
$ws.Range("M19").Value = 51.56497066666667
$ws.Range("N19").Value = 154.694912
$ws.Range("O19").Value = 0.2098600871655206
$ws.Range("P19").Value = 0.2187357342528896
$ws.Range("Q19").Value = 237.3489449137921
$ws.Range("R19").Value = 1424.093669482752
$ws.Range("S19").Value = 0.03437196548929499
$ws.Range("T19").Value = 0.02526301182107223
$ws.Range("G20").Value = 4.6029105
$ws.Range("H20").Value = 9.205821
$ws.Range("I20").Value = 0.1637851482553954
$ws.Range("J20").Value = 0.1154955860658076
$ws.Range("M20").Value = 29.9106925
$ws.Range("N20").Value = 59.821385
$ws.Range("O20").Value = 0.1217310987299521
$ws.Range("P20").Value = 0.08458632803643724
$ws.Range("Q20").Value = 137.6762405705213
$ws.Range("R20").Value = 550.704962282085
$ws.Range("S20").Value = 0.01993774605277738
$ws.Range("T20").Value = 0.009769347529722973
$ws.Range("G21").Value = 4.6029105
$ws.Range("H21").Value = 9.205821
$ws.Range("I21").Value = 0.1637851482553954
$ws.Range("J21").Value = 0.1154955860658076
$ws.Range("M21").Value = 70.603826
$ws.Range("N21").Value = 211.811478
$ws.Range("O21").Value = 0.2873447785777061
$ws.Range("P21").Value = 0.2994974984278718
$ws.Range("Q21").Value = 324.983092035573
$ws.Range("R21").Value = 1949.898552213438
$ws.Range("S21").Value = 0.04706280715976337
$ws.Range("T21").Value = 0.03459063910617034
$ws.Range("G22").Value = 1.953192
$ws.Range("H22").Value = 5.859576
$ws.Range("I22").Value = 0.0695003392508397
$ws.Range("J22").Value = 0.0735138304576138
$ws.Range("M22").Value = 44.04223000000001
$ws.Range("N22").Value = 132.12669
$ws.Range("O22").Value = 0.1792438957545786
$ws.Range("P22").Value = 0.1868246872369915
$ws.Range("Q22").Value = 86.02293129816002
$ws.Range("R22").Value = 774.2063816834401
$ws.Range("S22").Value = 0.01245751156358536
$ws.Range("T22").Value = 0.01373419838283692
$ws.Range("G23").Value = 1.953192
$ws.Range("H23").Value = 5.859576
$ws.Range("I23").Value = 0.0695003392508397
$ws.Range("J23").Value = 0.0735138304576138
$ws.Range("O23").Value = 0.2018201397722426
$ws.Range("P23").Value = 0.2103557520458098
$ws.Range("Q23").Value = 96.85774762440798
$ws.Range("R23").Value = 871.7197286196719
$ws.Range("S23").Value = 0.01402656818182274
$ws.Range("T23").Value = 0.01546405709167951
$ws.Range("G24").Value = 1.953192
$ws.Range("H24").Value = 5.859576
$ws.Range("I24").Value = 0.0695003392508397
$ws.Range("J24").Value = 0.0735138304576138
$ws.Range("M24").Value = 51.56497066666667
$ws.Range("N24").Value = 154.694912
$ws.Range("O24").Value = 0.2098600871655206
$ws.Range("P24").Value = 0.2187357342528896
$ws.Range("Q24").Value = 100.716288186368
$ws.Range("R24").Value = 906.4465936773121
$ws.Range("S24").Value = 0.01458534725321447
$ws.Range("T24").Value = 0.01608010168288859
$ws.Range("G25").Value = 1.953192
$ws.Range("H25").Value = 5.859576
$ws.Range("I25").Value = 0.0695003392508397
$ws.Range("J25").Value = 0.0735138304576138
$ws.Range("M25").Value = 29.9106925
$ws.Range("N25").Value = 59.821385
$ws.Range("O25").Value = 0.1217310987299521
$ws.Range("P25").Value = 0.08458632803643724
$ws.Range("Q25").Value = 58.42132530545999
$ws.Range("R25").Value = 350.52795183276
$ws.Range("S25").Value = 0.008460352659109129
$ws.Range("T25").Value = 0.006218264978302752
$ws.Range("G26").Value = 1.953192
$ws.Range("H26").Value = 5.859576
$ws.Range("I26").Value = 0.0695003392508397
$ws.Range("J26").Value = 0.0735138304576138
$ws.Range("M26").Value = 70.603826
$ws.Range("N26").Value = 211.811478
$ws.Range("O26").Value = 0.2873447785777061
$ws.Range("P26").Value = 0.2994974984278718
$ws.Range("Q26").Value = 137.902828112592
$ws.Range("R26").Value = 1241.125453013328
$ws.Range("S26").Value = 0.01997055959310799
$ws.Range("T26").Value = 0.02201720832190602
